# Khulo Municipality area sheet: drop the census-data footnote row and the
# 1989/2002 columns, keeping only the 2014 figure (export/maps fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held the "(according to the population census data)" footnote -
# remove it entirely so the rows below shift up one.
$ws.Rows("2:2").Delete()

# Columns B and C held the 1989 and 2002 figures - remove them so the
# 2014 column (old D) becomes the sheet's only data column (new B).
$ws.Columns("B:C").Delete()

# The new layout uses a taller, uniform row height for every row,
# including a few trailing blank rows left for future data.
$ws.Rows("1:8").RowHeight = 20.1
